$d = $word.ActiveDocument

for ($i = 1; $i -le 5; $i++) {
    $needle = "<id>p010r_$i</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assigning the exact same text as a no-op optimization can skip the
        # run merge, so first set a distinct placeholder, then set the
        # final text. This forces the three original runs (<id>, id value,
        # </id>) to collapse into a single run carrying the first run's
        # formatting (Courier New / 7f6000 / sz 18), matching the target.
        $rng.Text = $needle + "~"
        $rng.Text = $needle
    }
}
